# Katalog guncellendi - Cmt 29.11.2025 10:47:17,10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a typo in the image filename for row 121 (missing dot before extension)
$ws.Range("D121").Value = "SELANİKKOYUYEŞİL.jpg"

# Fix the "kategori" column (C) for the Kazak rows 118-124:
# header capitalization corrected from "KAZAK" to "Kazak"
for ($r = 118; $r -le 124; $r++) {
    $ws.Cells.Item($r, 3).Value = "Kazak"
}

# Reflect the final selection made by the user while editing
$ws.Range("C118:C124").Select()
